# Update LR-pair TPM-derived metrics in the NATMI output sheet.
# Corresponds to the commit "update scripts wuth new tpm": the ligand /
# receptor / edge expression values were recomputed with new TPM input,
# so this script overwrites the affected numeric cells (columns G-T,
# rows 2-7) with their newly recomputed values. Columns A-F are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs, Gnas -> Vipr1)
$ws.Range("G2").Value = 68.158272
$ws.Range("H2").Value = 204.474816
$ws.Range("I2").Value = 0.164824640128582
$ws.Range("J2").Value = 0.1648246401285819
$ws.Range("M2").Value = 0.8059226666666667
$ws.Range("N2").Value = 2.417768
$ws.Range("O2").Value = 0.1314814101815314
$ws.Range("P2").Value = 0.1314814101815314
$ws.Range("Q2").Value = 54.930296325632
$ws.Range("R2").Value = 494.372666930688
$ws.Range("S2").Value = 0.02167137611676938
$ws.Range("T2").Value = 0.02167137611676938

# Row 3 (ECs -> MuSCs)
$ws.Range("G3").Value = 68.158272
$ws.Range("H3").Value = 204.474816
$ws.Range("I3").Value = 0.164824640128582
$ws.Range("J3").Value = 0.1648246401285819
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.323633333333333
$ws.Range("N3").Value = 15.9709
$ws.Range("O3").Value = 0.8685185898184686
$ws.Range("P3").Value = 0.8685185898184687
$ws.Range("Q3").Value = 362.8496487616
$ws.Range("R3").Value = 3265.6468388544
$ws.Range("S3").Value = 0.1431532640118126
$ws.Range("T3").Value = 0.1431532640118126

# Row 4 (FAPs -> ECs)
$ws.Range("I4").Value = 0.3471155005059974
$ws.Range("J4").Value = 0.3471155005059974
$ws.Range("M4").Value = 0.8059226666666667
$ws.Range("N4").Value = 2.417768
$ws.Range("O4").Value = 0.1314814101815314
$ws.Range("P4").Value = 0.1314814101815314
$ws.Range("Q4").Value = 115.6814738812107
$ws.Range("R4").Value = 1041.133264930896
$ws.Range("S4").Value = 0.0456392355023966
$ws.Range("T4").Value = 0.0456392355023966

# Row 5 (FAPs -> MuSCs)
$ws.Range("I5").Value = 0.3471155005059974
$ws.Range("J5").Value = 0.3471155005059974
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.323633333333333
$ws.Range("N5").Value = 15.9709
$ws.Range("O5").Value = 0.8685185898184686
$ws.Range("P5").Value = 0.8685185898184687
$ws.Range("Q5").Value = 764.1499313455333
$ws.Range("R5").Value = 6877.3493821098
$ws.Range("S5").Value = 0.3014762650036008
$ws.Range("T5").Value = 0.3014762650036008

# Row 6 (MuSCs -> ECs)
$ws.Range("G6").Value = 201.822474
$ws.Range("H6").Value = 605.4674219999999
$ws.Range("I6").Value = 0.4880598593654206
$ws.Range("J6").Value = 0.4880598593654206
$ws.Range("M6").Value = 0.8059226666666667
$ws.Range("N6").Value = 2.417768
$ws.Range("O6").Value = 0.1314814101815314
$ws.Range("P6").Value = 0.1314814101815314
$ws.Range("Q6").Value = 162.653306439344
$ws.Range("R6").Value = 1463.879757954096
$ws.Range("S6").Value = 0.06417079856236538
$ws.Range("T6").Value = 0.06417079856236538

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 201.822474
$ws.Range("H7").Value = 605.4674219999999
$ws.Range("I7").Value = 0.4880598593654206
$ws.Range("J7").Value = 0.4880598593654206
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.323633333333333
$ws.Range("N7").Value = 15.9709
$ws.Range("O7").Value = 0.8685185898184686
$ws.Range("P7").Value = 0.8685185898184687
$ws.Range("Q7").Value = 1074.4288500022
$ws.Range("R7").Value = 9669.8596500198
$ws.Range("S7").Value = 0.4238890608030552
$ws.Range("T7").Value = 0.4238890608030553
